# Insert a new data row at row 373 (weekly price entry), shifting the
# existing rows 373-471 down to 374-472.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("373:373").Insert()

$ws.Range("A373").Value = 4
$ws.Range("B373").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C373").Value = "Los Lagos"
$ws.Range("D373").Value2 = 44943
$ws.Range("E373").Value = 10
$ws.Range("F373").Value = 100112008
$ws.Range("G373").Value = "Coliflor"
$ws.Range("H373").Value = "Sin especificar"
$ws.Range("I373").Value = "Primera"
$ws.Range("J373").Value = 1200
$ws.Range("K373").Value = 1600
$ws.Range("L373").Value = 1600
$ws.Range("M373").Value = 1600
$ws.Range("N373").Value = "$/unidad"
$ws.Range("O373").Value = "Región Metropolitana"
$ws.Range("P373").Value = 1600
$ws.Range("Q373").Value = 1
$ws.Range("R373").Value = "Hortaliza"
